# Add columns I (I0) and J (IF) to the worksheet, matching the style of
# the existing header row for the header cells, and plain numeric values
# for the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new header cells I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from an existing header cell (H1) so the new header
# cells match the bold/centered/bordered look of the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-set the values after paste (PasteSpecial formats only, so values are
# already intact, but ensure correctness).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-19
$values = @{
    2  = @(5, 5)
    3  = @(6, 7)
    4  = @(7, 7)
    5  = @(6, 7)
    6  = @(6, 7)
    7  = @(8, 8)
    8  = @(6, 7)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(7, 7)
    12 = @(6, 7)
    13 = @(6, 7)
    14 = @(7, 8)
    15 = @(7, 7)
    16 = @(5, 6)
    17 = @(8, 9)
    18 = @(5, 6)
    19 = @(9, 9)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($r, 10).Value = $pair[1]  # column J
}

$excel.CutCopyMode = 0
